# Insert a new weekly record at row 76 (Jengibre, Mercado Mayorista Lo Valledor
# de Santiago). This pushes the existing rows 76:86 down to 77:87, matching the
# target dimension A1:R87.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("76:76").Insert()

$ws.Cells.Item(76, 1).Value = 6
$ws.Cells.Item(76, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(76, 3).Value = "Metropolitana"
$ws.Cells.Item(76, 4).Value = 44748
$ws.Cells.Item(76, 5).Value = 13
$ws.Cells.Item(76, 6).Value = 100114007
$ws.Cells.Item(76, 7).Value = "Jengibre"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 220
$ws.Cells.Item(76, 11).Value = 11000
$ws.Cells.Item(76, 12).Value = 12000
$ws.Cells.Item(76, 13).Value = 11455
$ws.Cells.Item(76, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(76, 15).Value = "Perú"
$ws.Cells.Item(76, 16).Value = 881
$ws.Cells.Item(76, 17).Value = 13
$ws.Cells.Item(76, 18).Value = "Hortaliza"
